$wb = $excel.ActiveWorkbook

# xlPasteFormats = -4122
$xlPasteFormats = -4122

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

# C1 sits in the middle of the merged B1:D1 header box border -> keep only
# a top+bottom edge (the box border is being redrawn/split across the
# merged cells' underlying individual cells).
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.LineStyle = 1
$c1.Borders(7).LineStyle = 0
$c1.Borders(10).LineStyle = 0

# D1 sits at the right end of that box border -> keep top+right+bottom.
$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.LineStyle = 1
$d1.Borders(7).LineStyle = 0

# Rename header "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

# Sheet 2 has two merged header boxes (B1:D1 and E1:G1) with the same
# layout, so just copy the already-fixed formats from sheet 1's C1/D1
# instead of re-deriving them (avoids creating stray duplicate/unused
# cell styles for identical border combinations).
$c1.Copy()
$ws2.Range("C1").PasteSpecial($xlPasteFormats)
$ws2.Range("F1").PasteSpecial($xlPasteFormats)

$d1.Copy()
$ws2.Range("D1").PasteSpecial($xlPasteFormats)
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

# Rename headers "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell at G5
$ws2.Range("G5").ClearContents()
